$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 27 de Mayo de 2020 a las 22:35"

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 1740620
$ws.Range("C4").Value = 15345
$ws.Range("D4").Value = 483749
$ws.Range("E4").Value = 1155115
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 1184
$ws.Range("H4").Value = 101756

# Row 5: Brasil
$ws.Range("A5").Value = "Brasil"
$ws.Range("B5").Value = 399632
$ws.Range("C5").Value = 7272
$ws.Range("D5").Value = 158593
$ws.Range("E5").Value = 216004
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 486
$ws.Range("H5").Value = 25035

# Row 11: Alemania
$ws.Range("A11").Value = "Alemania"
$ws.Range("B11").Value = 181872
$ws.Range("C11").Value = 584
$ws.Range("D11").Value = 162800
$ws.Range("E11").Value = 10539
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 35
$ws.Range("H11").Value = 8533

# Row 33: Sudafrica
$ws.Range("A33").Value = "Sudafrica"
$ws.Range("B33").Value = 25937
$ws.Range("C33").Value = 1673
$ws.Range("D33").Value = 13451
$ws.Range("E33").Value = 11934
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 28
$ws.Range("H33").Value = 552

# Row 34: Irlanda
$ws.Range("A34").Value = "Irlanda"
$ws.Range("B34").Value = 24803
$ws.Range("C34").Value = 68
$ws.Range("D34").Value = 21060
$ws.Range("E34").Value = 2112
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 16
$ws.Range("H34").Value = 1631

# Row 96: Mayotte
$ws.Range("A96").Value = "Mayotte"
$ws.Range("B96").Value = 1645
$ws.Range("C96").Value = 11
$ws.Range("D96").Value = 1314
$ws.Range("E96").Value = 311
$ws.Range("F96").Value = 0
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 20

# Row 135: Congo
$ws.Range("A135").Value = "Congo"
$ws.Range("B135").Value = 571
$ws.Range("C135").Value = 84
$ws.Range("D135").Value = 161
$ws.Range("E135").Value = 391
$ws.Range("F135").Value = 0
$ws.Range("G135").Value = 3
$ws.Range("H135").Value = 19

# Row 136: Jamaica
$ws.Range("A136").Value = "Jamaica"
$ws.Range("B136").Value = 564
$ws.Range("C136").Value = 8
$ws.Range("D136").Value = 267
$ws.Range("E136").Value = 288
$ws.Range("F136").Value = 0
$ws.Range("G136").Value = 0
$ws.Range("H136").Value = 9

# Row 137: Tanzania
$ws.Range("A137").Value = "Tanzania"
$ws.Range("B137").Value = 509
$ws.Range("C137").Value = 0
$ws.Range("D137").Value = 183
$ws.Range("E137").Value = 305
$ws.Range("F137").Value = 0
$ws.Range("G137").Value = 0
$ws.Range("H137").Value = 21

# Row 142: Togo
$ws.Range("A142").Value = "Togo"
$ws.Range("B142").Value = 395
$ws.Range("C142").Value = 4
$ws.Range("D142").Value = 183
$ws.Range("E142").Value = 199
$ws.Range("F142").Value = 0
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 13

# Row 150: Mauritania
$ws.Range("A150").Value = "Mauritania"
$ws.Range("B150").Value = 292
$ws.Range("C150").Value = 24
$ws.Range("D150").Value = 15
$ws.Range("E150").Value = 261
$ws.Range("F150").Value = 0
$ws.Range("G150").Value = 3
$ws.Range("H150").Value = 16

# Row 151: Suazilandia
$ws.Range("A151").Value = "Suazilandia"
$ws.Range("B151").Value = 272
$ws.Range("C151").Value = 11
$ws.Range("D151").Value = 168
$ws.Range("E151").Value = 102
$ws.Range("F151").Value = 0
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 2

# Row 175: Barbados
$ws.Range("A175").Value = "Barbados"
$ws.Range("B175").Value = 92
$ws.Range("C175").Value = 0
$ws.Range("D175").Value = 76
$ws.Range("E175").Value = 9
$ws.Range("F175").Value = 0
$ws.Range("G175").Value = 0
$ws.Range("H175").Value = 7
